$d = $word.ActiveDocument

# 1. Remove the whole "Absolutely, I can analyze..." paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Absolutely, I can analyze the provided MATLAB code*") {
        $p.Range.Delete()
        break
    }
}

# 2. Trim the trailing parenthetical from the "It returns the average MSE..." sentence.
$d.Content.Find.Execute(
    "It returns the average MSE of the training, validation, and testing sets. (This might be a modification to the original code for potentially better performance based on the comments)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It returns the average MSE of the training, validation, and testing sets.", 2)

# 3. Drop "likely" from the transposed-input explanation.
$d.Content.Find.Execute(
    ". This is likely because the neural network expects the input data in transposed form (column vectors).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". This is because the neural network expects the input data in transposed form (column vectors).", 2)
